$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two now-obsolete trailing rows (52 and 53); Excel will
# automatically shrink the used range / merged cells (e.g. A8:O53 -> A8:O51).
$ws.Range("A52:A53").EntireRow.Delete()

# Apply the updated "PARTS" descriptions, confidence-qualifier fluid/material
# labels, and the stripped " G" gauge suffixes on design/operating pressure.
$ws.Range("D8").Value = "Thermal Deaerator"
$ws.Range("E8").Value = "Shell Plate"
$ws.Range("G8").Value = "Hot Water"
$ws.Range("M8").Value = "4 Bar"
$ws.Range("O8").Value = "1 Bar"
$ws.Range("E9").Value = "2:1 Ellipsoidal Head"
$ws.Range("G9").Value = "Hot Water"
$ws.Range("M9").Value = "4 Bar"
$ws.Range("O9").Value = "1 Bar"
$ws.Range("E10").Value = "Lifting Lug"
$ws.Range("G10").Value = "Hot Water"
$ws.Range("M10").Value = "4 Bar"
$ws.Range("O10").Value = "1 Bar"
$ws.Range("E11").Value = "Lifting Lug"
$ws.Range("G11").Value = "Hot Water"
$ws.Range("M11").Value = "4 Bar"
$ws.Range("O11").Value = "1 Bar"
$ws.Range("E12").Value = "Saddle Plate"
$ws.Range("G12").Value = "Hot Water"
$ws.Range("M12").Value = "4 Bar"
$ws.Range("O12").Value = "1 Bar"
$ws.Range("E13").Value = "Saddle Baseplate"
$ws.Range("G13").Value = "Hot Water"
$ws.Range("M13").Value = "4 Bar"
$ws.Range("O13").Value = "1 Bar"
$ws.Range("E14").Value = "Saddle Baseplate"
$ws.Range("G14").Value = "Hot Water"
$ws.Range("M14").Value = "4 Bar"
$ws.Range("O14").Value = "1 Bar"
$ws.Range("E15").Value = "Saddle Baseplate"
$ws.Range("G15").Value = "Hot Water"
$ws.Range("M15").Value = "4 Bar"
$ws.Range("O15").Value = "1 Bar"
$ws.Range("E16").Value = "Flange"
$ws.Range("G16").Value = "Hot Water"
$ws.Range("M16").Value = "4 Bar"
$ws.Range("O16").Value = "1 Bar"
$ws.Range("E17").Value = "Flange"
$ws.Range("G17").Value = "Hot Water"
$ws.Range("M17").Value = "4 Bar"
$ws.Range("O17").Value = "1 Bar"
$ws.Range("E18").Value = "Neck"
$ws.Range("G18").Value = "Hot Water"
$ws.Range("M18").Value = "4 Bar"
$ws.Range("O18").Value = "1 Bar"
$ws.Range("E19").Value = "Neck"
$ws.Range("G19").Value = "Hot Water"
$ws.Range("M19").Value = "4 Bar"
$ws.Range("O19").Value = "1 Bar"
$ws.Range("E20").Value = "Flange"
$ws.Range("G20").Value = "Hot Water"
$ws.Range("M20").Value = "4 Bar"
$ws.Range("O20").Value = "1 Bar"
$ws.Range("E21").Value = "Flange"
$ws.Range("G21").Value = "Hot Water"
$ws.Range("M21").Value = "4 Bar"
$ws.Range("O21").Value = "1 Bar"
$ws.Range("E22").Value = "Neck"
$ws.Range("G22").Value = "Hot Water"
$ws.Range("M22").Value = "4 Bar"
$ws.Range("O22").Value = "1 Bar"
$ws.Range("E23").Value = "Flange"
$ws.Range("G23").Value = "Hot Water"
$ws.Range("M23").Value = "4 Bar"
$ws.Range("O23").Value = "1 Bar"
$ws.Range("E24").Value = "Neck"
$ws.Range("G24").Value = "Hot Water"
$ws.Range("M24").Value = "4 Bar"
$ws.Range("O24").Value = "1 Bar"
$ws.Range("E25").Value = "Neck"
$ws.Range("G25").Value = "Hot Water"
$ws.Range("M25").Value = "4 Bar"
$ws.Range("O25").Value = "1 Bar"
$ws.Range("E26").Value = "Flange"
$ws.Range("G26").Value = "Hot Water"
$ws.Range("M26").Value = "4 Bar"
$ws.Range("O26").Value = "1 Bar"
$ws.Range("E27").Value = "Neck"
$ws.Range("G27").Value = "Hot Water"
$ws.Range("I27").Value = "ASTM A312"
$ws.Range("J27").Value = "TP304L"
$ws.Range("M27").Value = "4 Bar"
$ws.Range("O27").Value = "1 Bar"
$ws.Range("E28").Value = "Neck"
$ws.Range("G28").Value = "Hot Water"
$ws.Range("M28").Value = "4 Bar"
$ws.Range("O28").Value = "1 Bar"
$ws.Range("E29").Value = "Blind Flange"
$ws.Range("G29").Value = "Hot Water"
$ws.Range("I29").Value = "ASTM A182"
$ws.Range("J29").Value = "F304L"
$ws.Range("M29").Value = "4 Bar"
$ws.Range("O29").Value = "1 Bar"
$ws.Range("E30").Value = "Spiral Wound Gasket"
$ws.Range("G30").Value = "Hot Water"
$ws.Range("H30").Value = "Not Found"
$ws.Range("I30").Value = "ASME B16.20"
$ws.Range("J30").Value = "DN600"
$ws.Range("M30").Value = "4 Bar"
$ws.Range("O30").Value = "1 Bar"
$ws.Range("E31").Value = "Flange"
$ws.Range("G31").Value = "Hot Water"
$ws.Range("H31").Value = "Stainless Steel"
$ws.Range("I31").Value = "ASTM A182"
$ws.Range("J31").Value = "F304L"
$ws.Range("M31").Value = "4 Bar"
$ws.Range("O31").Value = "1 Bar"
$ws.Range("E32").Value = "Neck"
$ws.Range("G32").Value = "Hot Water"
$ws.Range("I32").Value = "ASTM A240"
$ws.Range("J32").Value = "304L"
$ws.Range("M32").Value = "4 Bar"
$ws.Range("O32").Value = "1 Bar"
$ws.Range("E33").Value = "Stud Bolt"
$ws.Range("G33").Value = "Hot Water"
$ws.Range("H33").Value = "Stainless Steel Bolting"
$ws.Range("I33").Value = "ASTM A193"
$ws.Range("J33").Value = "GR B8M"
$ws.Range("M33").Value = "4 Bar"
$ws.Range("O33").Value = "1 Bar"
$ws.Range("E34").Value = "Nuts & Washer"
$ws.Range("G34").Value = "Hot Water"
$ws.Range("H34").Value = "Heavy Hex Nuts"
$ws.Range("I34").Value = "ASTM A194"
$ws.Range("J34").Value = "GR 2H"
$ws.Range("M34").Value = "4 Bar"
$ws.Range("O34").Value = "1 Bar"
$ws.Range("E35").Value = "Bracket 1"
$ws.Range("G35").Value = "Hot Water"
$ws.Range("H35").Value = "Not Found"
$ws.Range("I35").Value = "ASTM A36"
$ws.Range("J35").Value = "-"
$ws.Range("M35").Value = "4 Bar"
$ws.Range("O35").Value = "1 Bar"
$ws.Range("E36").Value = "Angle Bar"
$ws.Range("G36").Value = "Hot Water"
$ws.Range("J36").Value = "-"
$ws.Range("M36").Value = "4 Bar"
$ws.Range("O36").Value = "1 Bar"
$ws.Range("E37").Value = "Bracket 2"
$ws.Range("G37").Value = "Hot Water"
$ws.Range("J37").Value = "-"
$ws.Range("M37").Value = "4 Bar"
$ws.Range("O37").Value = "1 Bar"
$ws.Range("E38").Value = "Bracket 3"
$ws.Range("G38").Value = "Hot Water"
$ws.Range("J38").Value = "-"
$ws.Range("M38").Value = "4 Bar"
$ws.Range("O38").Value = "1 Bar"
$ws.Range("E39").Value = "Grating"
$ws.Range("G39").Value = "Hot Water"
$ws.Range("J39").Value = "-"
$ws.Range("M39").Value = "4 Bar"
$ws.Range("O39").Value = "1 Bar"
$ws.Range("E40").Value = "Angle Bar"
$ws.Range("G40").Value = "Hot Water"
$ws.Range("J40").Value = "-"
$ws.Range("M40").Value = "4 Bar"
$ws.Range("O40").Value = "1 Bar"
$ws.Range("E41").Value = "Angle Bar"
$ws.Range("G41").Value = "Hot Water"
$ws.Range("J41").Value = "-"
$ws.Range("M41").Value = "4 Bar"
$ws.Range("O41").Value = "1 Bar"
$ws.Range("E42").Value = "Flat Plate"
$ws.Range("G42").Value = "Hot Water"
$ws.Range("J42").Value = "-"
$ws.Range("M42").Value = "4 Bar"
$ws.Range("O42").Value = "1 Bar"
$ws.Range("E43").Value = "Flat Plate"
$ws.Range("G43").Value = "Hot Water"
$ws.Range("J43").Value = "-"
$ws.Range("M43").Value = "4 Bar"
$ws.Range("O43").Value = "1 Bar"
$ws.Range("E44").Value = "Angle Bar"
$ws.Range("G44").Value = "Hot Water"
$ws.Range("J44").Value = "-"
$ws.Range("M44").Value = "4 Bar"
$ws.Range("O44").Value = "1 Bar"
$ws.Range("E45").Value = "Angle Bar"
$ws.Range("G45").Value = "Hot Water"
$ws.Range("J45").Value = "-"
$ws.Range("M45").Value = "4 Bar"
$ws.Range("O45").Value = "1 Bar"
$ws.Range("E46").Value = "Side Rail"
$ws.Range("G46").Value = "Hot Water"
$ws.Range("J46").Value = "-"
$ws.Range("M46").Value = "4 Bar"
$ws.Range("O46").Value = "1 Bar"
$ws.Range("E47").Value = "Ladder Bracket"
$ws.Range("G47").Value = "Hot Water"
$ws.Range("J47").Value = "-"
$ws.Range("M47").Value = "4 Bar"
$ws.Range("O47").Value = "1 Bar"
$ws.Range("E48").Value = "Round Bar"
$ws.Range("G48").Value = "Hot Water"
$ws.Range("J48").Value = "-"
$ws.Range("M48").Value = "4 Bar"
$ws.Range("O48").Value = "1 Bar"
$ws.Range("E49").Value = "Saddle Support Plate"
$ws.Range("G49").Value = "Hot Water"
$ws.Range("H49").Value = "Stainless Steel"
$ws.Range("I49").Value = "ASTM A240"
$ws.Range("J49").Value = "J304L"
$ws.Range("M49").Value = "4 Bar"
$ws.Range("O49").Value = "1 Bar"
$ws.Range("E50").Value = "Pad Plate"
$ws.Range("G50").Value = "Hot Water"
$ws.Range("J50").Value = "-"
$ws.Range("M50").Value = "4 Bar"
$ws.Range("O50").Value = "1 Bar"
$ws.Range("E51").Value = "Pad Plate"
$ws.Range("G51").Value = "Hot Water"
$ws.Range("H51").Value = "Not Found"
$ws.Range("I51").Value = "ASTM A36"
$ws.Range("J51").Value = "-"
$ws.Range("M51").Value = "4 Bar"
$ws.Range("O51").Value = "1 Bar"
